$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# New order rows (Remessa / Material / Quantidade) to populate rows 94-102,
# which previously existed as blank placeholder rows.
$data = @(
    @(94,  "80266387", "21012-CTY-I", 2),
    @(95,  "80266387", "20941-CTY-I", 5),
    @(96,  "80266387", "21016-CTY-I", 3),
    @(97,  "80266388", "21012-CTY-I", 1),
    @(98,  "80266388", "21013-CTY-I", 3),
    @(99,  "80266388", "21014-CTY-I", 5),
    @(100, "80266389", "10253-ARI-I", 1),
    @(101, "80266391", "10382-ARI-I", 1),
    @(102, "80266392", "10638-ARI-I", 1)
)

foreach ($row in $data) {
    $r = $row[0]

    # Column A ("Remessa") values are purely numeric-looking codes that must be
    # stored as text (matching the rest of the column). Writing them with
    # .Value directly would make Excel coerce them to numbers, so instead
    # build them as a text formula first, then convert that formula to a
    # static value in place - this keeps the original "General" cell style
    # untouched while still persisting the cell as text.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Formula = '="' + $row[1] + '"'
    $cellA.Copy() | Out-Null
    $cellA.PasteSpecial(-4163) | Out-Null  # xlPasteValues

    # Column B ("Material") values already contain letters/dashes so they are
    # stored as text automatically.
    $ws.Cells.Item($r, 2).Value = $row[2]

    # Column C ("Quantidade") is a plain number.
    $ws.Cells.Item($r, 3).Value = $row[3]
}

$excel.CutCopyMode = 0

# Update the selected range on the active sheet view to reflect the new data extent.
$ws.Range("A2:C102").Select()
